$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# NOTE: this runtime re-coalesces adjacent runs that share identical rPr
# whenever *any* run inside their paragraph is mutated (Find/Replace or
# Range.Text assignment). That is fine when the two runs touching an edit
# already differ in formatting (bold heading vs. plain body text, as in most
# paragraphs below) but it silently glues together unrelated, untouched runs
# when they happen to share formatting (see Change 2, which has 4 runs with
# matching rPr on runs 2-4). To keep those run boundaries intact we briefly
# flip Italic on the runs bordering the ones we edit (breaking the "same
# formatting" identity with their neighbour for the duration of the edit)
# and flip it back once the text is in place.
# ---------------------------------------------------------------------------

# --- Change 1 ---------------------------------------------------------------
# "Limited Brand Recognition and Awareness" (bold run) + following (non-bold)
# sentence, both reworded. The two runs already differ in boldness so no
# extra care is required to keep them apart.
$d.Content.Find.Execute("Limited Brand Recognition and Awareness", $true, $false, $false, $false, $false, $true, 1, $false, "Begrenzte Markenbekanntheit und Awareness", 2)

$d.Content.Find.Execute(": Die Sichtbarkeit in diesen neuen Märkten zu erreichen, ist eine primäre Hürde, die robuste Marketingbemühungen erfordert, um die Markenpräsenz von Adatum von Grund auf aufzubauen.", $true, $false, $false, $false, $false, $true, 1, $false, ": Das Erreichen von Sichtbarkeit in diesen neuen Märkten ist eine der größten Hürden und erfordert starke Marketingaktivitäten, um die Markenpräsenz von Adatum von Grund auf aufzubauen.", 2)

# --- Change 2 -----------------------------------------------------------------
# Paragraph layout: [bold "Intensiver Wettbewerb"] [run: "Der Cloud..." ] [run: " "] [run: "Adatum muss..."]
# The last three runs share identical rPr, so editing the first of them would
# otherwise cause the engine to glue all three back together. Break the
# identity around the edited run and the trailing run, edit, then restore.
$rTarget = $d.Content
$rTarget.Find.Execute(": Der Cloud services-Sektor in Kanada ist mit zahlreichen Akteuren stark wettbewerbsfähig.", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$rTarget.Italic = 1
$rTarget.Text = ": Der kanadische Cloud-Services-Sektor ist wettbewerbsintensiv und hat viele Akteure."
$tStart = $rTarget.Start
$tEnd = $rTarget.End

$rAfter = $d.Content
$rAfter.Find.Execute("Adatum muss den einzigartigen Wert seiner Lösungen klar artikulieren, um sich eine Nische zu schaffen.", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$rAfter.Italic = 1
$afterStart = $rAfter.Start
$afterEnd = $rAfter.End

$d.Range($tStart, $tEnd).Italic = 0
$d.Range($afterStart, $afterEnd).Italic = 0

# --- Change 3 -------------------------------------------------------------
$d.Content.Find.Execute(": Das Anpassen von Produkten und Marketing, um den unterschiedlichen Anforderungen dieser Märkte gerecht zu werden, ist entscheidend für das Resonieren mit lokalen Unternehmen und Verbrauchern.", $true, $false, $false, $false, $false, $true, 1, $false, ": Das Anpassen von Produkten und Marketing, um den unterschiedlichen Anforderungen dieser Märkte gerecht zu werden, ist entscheidend, um bei lokalen Unternehmen und Verbrauchern gut anzukommen.", 2)

# --- Change 4 ---------------------------------------------------------------
# The bold/non-bold run boundary shifts: "gesetzlichen Vorschriften und
# Compliance" moves from the (non-bold) second run into the (bold) heading
# run. Rewrite the heading run's Range.Text in place (keeps the edit inside
# the same bold run, no run split) then fix up the remainder of the sentence
# (still in the trailing non-bold run) via Find/Replace.
$rHeading = $d.Content
$rHeading.Find.Execute("Herausforderungen bei", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$rHeading.Text = "Herausforderungen bei gesetzlichen Vorschriften und Compliance"

$d.Content.Find.Execute(" gesetzlichen Vorschriften und Compliance: Adatum steht vor der komplexen Aufgabe, in den unterschiedlichen Datenschutz-, Sicherheits- und Betrieblichen Vorschriften der Region zu navigieren, was eine sorgfältige Einhaltung der Vorschriften erfordert.", $true, $false, $false, $false, $false, $true, 1, $false, ": Adatum steht vor der komplexen Aufgabe, sich in den verschiedenen Datenschutz-, Sicherheits- und Betriebsvorschriften der Region zurechtzufinden, was eine sorgfältige Einhaltung der Vorschriften erfordert.", 2)

# --- Change 5 -------------------------------------------------------------
$d.Content.Find.Execute("Operative und logistische Komplexitäten", $true, $false, $false, $false, $false, $true, 1, $false, "Operative und logistische Komplexität", 2)
